$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "70.098.22"
$ws.Range("E2").Value = "  +1.15%  "
$ws.Range("D3").Value = "3.504.73"
$ws.Range("E3").Value = "  +0.30%  "
$ws.Range("E4").Value = "  -0.07%  "
$ws.Range("D5").Value = "'602.40"
$ws.Range("E5").Value = "  -0.37%  "
$ws.Range("D6").Value = "'174.90"
$ws.Range("E6").Value = "  +3.70%  "
$ws.Range("E7").Value = "  -1.31%  "
$ws.Range("D8").Value = "3.499.60"
$ws.Range("E8").Value = "  +0.22%  "
$ws.Range("E9").Value = "  +0.02%  "
$ws.Range("E10").Value = "  +0.31%  "
$ws.Range("D11").Value = "'7.22"
$ws.Range("E11").Value = "  +9.15%  "
$ws.Range("E12").Value = "  +0.82%  "
$ws.Range("D13").Value = "'46.15"
$ws.Range("E13").Value = "  -1.27%  "
$ws.Range("D14").Value = "'0.0000275"
$ws.Range("E14").Value = "  -0.41%  "
$ws.Range("D15").Value = "4.071.49"
$ws.Range("E15").Value = "  +0.36%  "
$ws.Range("E16").Value = "  +0.14%  "
$ws.Range("D17").Value = "'610.30"
$ws.Range("E17").Value = "  +0.35%  "
$ws.Range("D18").Value = "3.504.16"
$ws.Range("E18").Value = "  +0.11%  "
$ws.Range("D19").Value = "70.192.57"
$ws.Range("E20").Value = "  +1.09%  "
$ws.Range("D21").Value = "'17.33"
$ws.Range("E21").Value = "  +1.11%  "
$ws.Range("E22").Value = "  +0.12%  "
$ws.Range("D23").Value = "'8.99"
$ws.Range("E23").Value = "  -12.14%  "
$ws.Range("D24").Value = "'97.80"
$ws.Range("E24").Value = "  +2.58%  "
$ws.Range("D25").Value = "'15.51"
$ws.Range("E25").Value = "  -1.04%  "
$ws.Range("E26").Value = "  -3.34%  "
$ws.Range("E27").Value = "  +0.05%  "
$ws.Range("D28").Value = "'2.55"
$ws.Range("E28").Value = "  -1.44%  "
$ws.Range("D29").Value = "'33.77"
$ws.Range("E29").Value = "  +2.55%  "
$ws.Range("E30").Value = "  -2.30%  "
$ws.Range("D31").Value = "'2.96"
$ws.Range("E31").Value = "  -3.09%  "
$ws.Range("D32").Value = "'8.00"
$ws.Range("E32").Value = "  -4.42%  "
$ws.Range("D33").Value = "'636.50"
$ws.Range("E33").Value = "  +15.15%  "
$ws.Range("E34").Value = "  -3.47%  "
$ws.Range("D35").Value = "'6.83"
$ws.Range("E35").Value = "  +0.06%  "
$ws.Range("B36").Value = "Hedera"
$ws.Range("C36").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D36").Value = "'0.0991"
$ws.Range("E36").Value = "  -1.33%  "
$ws.Range("B37").Value = "dogwifhat"
$ws.Range("C37").Value = "https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif"
$ws.Range("D37").Value = "'3.54"
$ws.Range("E37").Value = "  +2.66%  "
$ws.Range("D38").Value = "'10.72"
$ws.Range("E38").Value = "  +0.14%  "
$ws.Range("D39").Value = "'0.0471"
$ws.Range("E39").Value = "  +5.68%  "
$ws.Range("D40").Value = "'56.64"
$ws.Range("E40").Value = "  +0.13%  "
$ws.Range("E41").Value = "  -0.04%  "
$ws.Range("E42").Value = "  +2.51%  "
$ws.Range("D43").Value = "3.352.32"
$ws.Range("E43").Value = "  +0.38%  "
$ws.Range("D44").Value = "0.0₃0735"
$ws.Range("E44").Value = "  +6.41%  "
$ws.Range("E45").Value = "  -4.82%  "
$ws.Range("D46").Value = "'32.10"
$ws.Range("E46").Value = "  -2.05%  "
$ws.Range("D47").Value = "'2.88"
$ws.Range("E47").Value = "  +0.69%  "
$ws.Range("E48").Value = "  -2.33%  "
$ws.Range("E49").Value = "  +0.38%  "
$ws.Range("D50").Value = "'132.87"
$ws.Range("E50").Value = "  -1.16%  "
